$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/artfynd/A 30779-2023.xlsx")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/kartor/A 30779-2023.png")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/knärot/A 30779-2023.png")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/klagomål/A 30779-2023.docx")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/klagomålsmail/A 30779-2023.docx")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/tillsyn/A 30779-2023.docx")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/tillsynsmail/A 30779-2023.docx")'
